$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$sub3 = [char]0x2083

Set-TextValue $ws "D2" "65.451.58"
Set-TextValue $ws "E2" "  -5.81%  "

Set-TextValue $ws "D3" "3.295.04"
Set-TextValue $ws "E3" "  -7.00%  "

Set-TextValue $ws "D4" "0.999"
Set-TextValue $ws "E4" "  -0.15%  "

Set-TextValue $ws "D5" "554.40"
Set-TextValue $ws "E5" "  -5.23%  "

Set-TextValue $ws "D6" "178.87"
Set-TextValue $ws "E6" "  -8.53%  "

Set-TextValue $ws "E7" "  +0.08%  "

Set-TextValue $ws "D8" "0.584"
Set-TextValue $ws "E8" "  -4.37%  "

Set-TextValue $ws "D9" "3.287.97"
Set-TextValue $ws "E9" "  -6.87%  "

Set-TextValue $ws "D10" "0.182"
Set-TextValue $ws "E10" "  -11.32%  "

Set-TextValue $ws "D11" "0.579"
Set-TextValue $ws "E11" "  -8.09%  "

Set-TextValue $ws "D12" "46.78"
Set-TextValue $ws "E12" "  -11.33%  "

Set-TextValue $ws "D13" "0.0000260"
Set-TextValue $ws "E13" "  -9.51%  "

Set-TextValue $ws "D14" "3.824.72"
Set-TextValue $ws "E14" "  -7.06%  "

Set-TextValue $ws "D15" "8.45"
Set-TextValue $ws "E15" "  -8.39%  "

Set-TextValue $ws "D16" "594.03"
Set-TextValue $ws "E16" "  -10.57%  "

Set-TextValue $ws "D17" "65.461.21"
Set-TextValue $ws "E17" "  -5.99%  "

Set-TextValue $ws "D18" "17.86"
Set-TextValue $ws "E18" "  -2.97%  "

Set-TextValue $ws "D19" "0.117"
Set-TextValue $ws "E19" "  -4.09%  "

Set-TextValue $ws "D20" "3.292.06"
Set-TextValue $ws "E20" "  -7.44%  "

Set-TextValue $ws "D21" "11.29"
Set-TextValue $ws "E21" "  -9.45%  "

Set-TextValue $ws "D22" "0.890"
Set-TextValue $ws "E22" "  -7.56%  "

Set-TextValue $ws "D23" "16.80"
Set-TextValue $ws "E23" "  -6.96%  "

Set-TextValue $ws "D24" "5.03"
Set-TextValue $ws "E24" "  -6.46%  "

Set-TextValue $ws "D25" "100.01"
Set-TextValue $ws "E25" "  -4.68%  "

Set-TextValue $ws "D26" "3.95"
Set-TextValue $ws "E26" "  -9.85%  "

Set-TextValue $ws "D27" "5.99"
Set-TextValue $ws "E27" "  -0.32%  "

Set-TextValue $ws "D28" "2.64"
Set-TextValue $ws "E28" "  -9.93%  "

Set-TextValue $ws "D29" "9.16"
Set-TextValue $ws "E29" "  -9.61%  "

Set-TextValue $ws "D30" "8.57"
Set-TextValue $ws "E30" "  -10.41%  "

Set-TextValue $ws "D31" "30.26"
Set-TextValue $ws "E31" "  -9.04%  "

Set-TextValue $ws "D32" "3.81"
Set-TextValue $ws "E32" "  -11.29%  "

Set-TextValue $ws "D33" "6.18"
Set-TextValue $ws "E33" "  -8.50%  "

Set-TextValue $ws "D34" "10.90"
Set-TextValue $ws "E34" "  -7.22%  "

Set-TextValue $ws "E35" "  -7.39%  "

$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D36" "3.753.73"
Set-TextValue $ws "E36" "  -0.70%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D37" "57.29"
Set-TextValue $ws "E37" "  -7.60%  "

Set-TextValue $ws "D38" "1.00"
Set-TextValue $ws "E38" "  +0.16%  "

Set-TextValue $ws "D39" "510.90"
Set-TextValue $ws "E39" "  +2.43%  "

Set-TextValue $ws "D40" "3.44"
Set-TextValue $ws "E40" "  -8.87%  "

Set-TextValue $ws "D41" ("0.0{0}0703" -f $sub3)
Set-TextValue $ws "E41" "  -12.55%  "

Set-TextValue $ws "D42" "2.62"
Set-TextValue $ws "E42" "  -8.80%  "

Set-TextValue $ws "D43" "0.124"
Set-TextValue $ws "E43" "  -8.03%  "

Set-TextValue $ws "D44" "0.335"
Set-TextValue $ws "E44" "  -9.37%  "

Set-TextValue $ws "D45" "31.55"
Set-TextValue $ws "E45" "  -8.99%  "

Set-TextValue $ws "D46" "3.25"
Set-TextValue $ws "E46" "  -3.94%  "

Set-TextValue $ws "D47" "0.0408"
Set-TextValue $ws "E47" "  -9.35%  "

Set-TextValue $ws "D48" "3.05"
Set-TextValue $ws "E48" "  +13.81%  "

Set-TextValue $ws "D49" "0.128"
Set-TextValue $ws "E49" "  -5.93%  "

Set-TextValue $ws "E50" "  -10.61%  "

Set-TextValue $ws "D51" "0.997"
Set-TextValue $ws "E51" "  -0.49%  "
